$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A3:F11")
Write-Output ("before color: " + $r.Font.Color)
$r.HorizontalAlignment = -4131
Write-Output ("after color: " + $r.Font.Color)
Write-Output "done"
